$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Column D values are numeric-looking strings that must
# remain TEXT (matching the original t="inlineStr" cells), so we prefix them with a
# leading apostrophe - same as typing '1.00 into Excel - which forces text storage
# (quotePrefix) without permanently changing the cell NumberFormat.

$ws.Range("D2").Value = "'68.506.18"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "'3.773.30"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'595.88"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'168.63"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'3.770.12"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("D11").Value = "'6.49"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "'36.59"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'4.413.44"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'3.784.06"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'68.495.79"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'18.25"
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("D19").Value = "'7.07"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("D22").Value = "'469.53"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'0.703"
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("D24").Value = "'84.49"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'0.0000143"
$ws.Range("E25").Value = "  -4.17%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'12.28"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "'10.26"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'3.928.37"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'2.80"
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("D32").Value = "'7.45"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'30.18"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("D37").Value = "'3.734.14"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("E39").Value = "  -8.57%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "'5.83"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.308"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "'1.98"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "'43.66"
$ws.Range("E47").Value = "  +12.69%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").Value = "'408.20"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'45.47"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "'145.16"
$ws.Range("E51").Value = "  +2.17%  "
